# aisetting.xlsx — add a 6th prompt column ("a6") to the settings sheet.
#
# Layout before the edit:
#   Row 1 (headers): A1="a1" B1="a2" C1="a3" D1="a4" E1="a5"
#   Row 2..5 hold the per-column prompt text (col E uses a wrapped, very
#   tall cell, style index 1 == WrapText).
#   E5 currently holds the long "Provide a list of requested number of
#   items. ..." prompt.
#
# After the edit:
#   A new column F is introduced ("a6" header), the big "Provide a list of
#   requested number of items. ..." prompt moves from E5 down into the new
#   F2 cell, and E5 gets a brand new prompt text. The view scroll/selection
#   is reset (selection -> E6, no frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the long prompt currently sitting in E5 down into the new F2 cell
# (Cut+paste keeps it as the same underlying string entry instead of
# minting a duplicate).
$ws.Range("E5").Cut($ws.Range("F2"))
$ws.Range("F2").WrapText = $true

# --- New column F header + width ---------------------------------------
$ws.Range("F1").Value = "a6"
$ws.Columns.Item(6).ColumnWidth = 21.72

# --- Updated E5 prompt --------------------------------------------------
$newPrompt = "Provide a list of 20 most related best questions with answers, in this format: `r`nitem x/{question}/{answer}.`r`nFinal output are in the following format:`r`n    - item 1`r`n    - item 2`r`n    - item 3"
$ws.Range("E5").Value = $newPrompt
$ws.Range("E5").WrapText = $true

# --- Reset the view: selection on E6, scrolled back to the top-left ----
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$null = $ws.Range("E6").Select()

Write-Output "done"
